# Stepper BOM final update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update connector part numbers (Comment column A) for J1, J2 and J4
$ws.Range("A18").Value = "5-146280-2"
$ws.Range("A19").Value = "5-146280-4"
$ws.Range("A21").Value = "XH2.54MM 4pins right"

# Restore the view/selection as left by the author (scrolled back to top,
# with A22 selected instead of A24)
$ws.Range("A22").Select()
